$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2" = 1.05
    "H2" = 1.05
    "J2" = 1.06
    "K2" = 980
    "N2" = 1.3
    "O2" = 1.24
    "Q2" = 1.24
    "S2" = 1.05
    "T2" = 1.04
    "U2" = 1.04
    "V2" = 1.02
    "W2" = 1.02
    "K3" = 4.3
    "V3" = 1.23
    "W4" = 1.2
    "H6" = 1.85
    "F8" = 1.97
    "V8" = 1.3
    "S9" = 3.1
    "G10" = 5.1
    "H10" = 1.85
    "J10" = 3.4
    "L11" = 1.23
    "S11" = 2.2
    "T11" = 2.04
    "U11" = 1.77
    "AD11" = 55
    "J12" = 3.5
    "N13" = 1.57
    "P13" = 1.57
    "F14" = 2.28
    "H14" = 2.76
    "N14" = 5.2
    "P14" = 2.42
    "R14" = 1.58
    "S14" = 2.38
    "AN14" = 15.5
    "G15" = 11
    "H15" = 1.41
    "I15" = 1.42
    "J15" = 5.1
    "K15" = 5.3
    "O15" = 1.33
    "V15" = 3.35
    "Y15" = 7.4
    "AF15" = 100
    "AJ15" = 520
    "AK15" = 230
    "AN15" = 360
    "T16" = 1.52
    "W16" = 1.25
    "F17" = 2.6
    "G17" = 2.9
    "H17" = 2.42
    "I17" = 2.68
    "K17" = 4.4
    "Q17" = 1.59
    "W17" = 1.52
    "AJ17" = 60
    "G18" = 1.84
    "J18" = 3.15
    "Q18" = 2.2
    "S18" = 3.85
    "U18" = 1.74
    "V18" = 1.15
    "W18" = 2.2
    "G19" = 2.6
    "L19" = 1.45
    "M19" = 1.08
    "N19" = 3.55
    "P19" = 1.84
    "R19" = 1.32
    "T19" = 1.86
    "AM19" = 100
    "F20" = 1.86
    "H20" = 5.3
    "J20" = 3.55
    "R20" = 1.24
    "T20" = 2.18
    "W20" = 2.14
    "AB20" = 7
    "AC20" = 7.8
    "F21" = 1.35
    "H21" = 11
    "K21" = 6
    "R21" = 1.47
    "W21" = 3.8
    "AC21" = 12.5
    "AE21" = 210
    "AI21" = 210
    "AM21" = 190
    "AN21" = 5.8
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
